$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2: Fecha del Reporte (date serial 45371 -> 45373)
$ws.Range("C2").Value = 45373

# C3: Ficha de Caracterizacion - must remain a text value ("2499992" -> "2824078")
# Temporarily mark the cell as Text so the numeric-looking string isn't
# auto-converted to a number, then restore the original (General) formatting.
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "2824078"
$ws.Range("C3").ClearFormats()

# C4: Programa (text)
$ws.Range("C4").Value = "DISEÑO E INTEGRACION DE MULTIMEDIA"

# C6: Fecha Inicio (date serial 44669 -> 44760)
$ws.Range("C6").Value = 44760

# C7: Fecha Fin (date serial 45490 -> 45124)
$ws.Range("C7").Value = 45124
